$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "67.943.11"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "3.743.64"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'594.80"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "'166.63"
$ws.Range("E6").Value = "  -1.38%  "
$ws.Range("D7").Value = "3.740.06"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").Value = "'6.48"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("E13").Value = "  -4.63%  "
$ws.Range("D14").Value = "'36.65"
$ws.Range("D15").Value = "4.373.23"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("D16").Value = "3.739.79"
$ws.Range("E16").Value = "  -1.68%  "
$ws.Range("D17").Value = "67.957.55"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'18.15"
$ws.Range("E18").Value = "  -3.19%  "
$ws.Range("E19").Value = "  -5.56%  "
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "'467.50"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("E23").Value = "  -4.97%  "
$ws.Range("D24").Value = "'83.17"
$ws.Range("E24").Value = "  -0.68%  "
$ws.Range("D25").Value = "'0.0000138"
$ws.Range("E25").Value = "  -8.89%  "
$ws.Range("E26").Value = "  -2.67%  "
$ws.Range("D27").Value = "'12.10"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("D30").Value = "3.892.52"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "'2.79"
$ws.Range("E31").Value = "  -4.70%  "
$ws.Range("D32").Value = "'7.39"
$ws.Range("E32").Value = "  -4.68%  "
$ws.Range("D33").Value = "'2.25"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").Value = "'29.86"
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D35").Value = "'9.13"
$ws.Range("E35").Value = "  -2.00%  "
$ws.Range("D36").Value = "'0.996"
$ws.Range("D37").Value = "3.697.28"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("E38").Value = "  -4.21%  "
$ws.Range("D39").Value = "'3.44"
$ws.Range("E39").Value = "  -10.98%  "
$ws.Range("D40").Value = "'0.137"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("E42").Value = "  -3.16%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("E46").Value = "  -2.01%  "
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").Value = "'45.34"
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("D49").Value = "'396.74"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("D50").Value = "'144.56"
$ws.Range("E50").Value = "  +1.15%  "
$ws.Range("E51").Value = "  +0.14%  "
